# Applies the commit "allocation version working fine":
#  - Renames "DOCs" -> "DOC_1" and replaces its distance-cost-matrix data with new values
#  - Moves "Aircraft" ahead of "DOC_1" in the tab order (keeping its name)
#  - Replaces "Aircraft" data: new min_load_factor column, refreshed payload/range strings,
#    and updated numbers
#  - Adds two new sheets "DOC_2" and "DOC_3" (further distance-cost matrices) after "DOC_1"
#  - Updates sheet view settings (zoom / selected cell) and the active tab ("Demands")

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Reorder: move "Aircraft" so it sits right before "DOCs" (soon "DOC_1")
# ---------------------------------------------------------------------------
$docs = $wb.Worksheets.Item("DOCs")
$aircraft = $wb.Worksheets.Item("Aircraft")
$aircraft.Move($docs)

# ---------------------------------------------------------------------------
# 2. Rewrite the "Aircraft" sheet contents
# ---------------------------------------------------------------------------
$aircraft = $wb.Worksheets.Item("Aircraft")
$aircraft.Cells.Clear()

$aircraft.Range("A1").Value = "acft_num"
$aircraft.Range("B1").Value = "passenger_weight"
$aircraft.Range("C1").Value = "min_load_factor"
$aircraft.Range("D1").Value = "avg_ticket_price"
$aircraft.Range("E1").Value = "max_allowed_acft"
$aircraft.Range("F1").Value = "payloads"
$aircraft.Range("G1").Value = "ranges"

$aircraft.Range("A2").Value = 1
$aircraft.Range("B2").Value = 100
$aircraft.Range("C2").Value = 0.5
$aircraft.Range("D2").Value = 120
$aircraft.Range("E2").Value = 10
$aircraft.Range("F2").Value = "[8000,8000,5000,0]"
$aircraft.Range("G2").Value = "[0,1000,2000,2500]"

$aircraft.Range("A3").Value = 1
$aircraft.Range("B3").Value = 100
$aircraft.Range("C3").Value = 0.5
$aircraft.Range("D3").Value = 110
$aircraft.Range("E3").Value = 10
$aircraft.Range("F3").Value = "[12000,12000,8000,0]"
$aircraft.Range("G3").Value = "[0,1000,2000,2500]"

$aircraft.Range("A4").Value = 1
$aircraft.Range("B4").Value = 100
$aircraft.Range("C4").Value = 0.5
$aircraft.Range("D4").Value = 100
$aircraft.Range("E4").Value = 10
$aircraft.Range("F4").Value = "[15000,15000,10000,0]"
$aircraft.Range("G4").Value = "[0,1000,2000,2500]"

# stray formatted (underlined) empty cell left a few rows below the table
$aircraft.Range("G10").Font.Underline = 1

$aircraft.Activate()
$wb.Windows.Item(1).Zoom = 160
$aircraft.Range("C3").Select()

# ---------------------------------------------------------------------------
# 3. Rename "DOCs" -> "DOC_1" and rewrite its distance-cost matrix
# ---------------------------------------------------------------------------
$docs = $wb.Worksheets.Item("DOCs")
$docs.Name = "DOC_1"

$doc1 = $wb.Worksheets.Item("DOC_1")

$doc1.Range("A1").Value = 0
$doc1.Range("B1").Value = 4781.7
$doc1.Range("C1").Value = 3394.2999999999997
$doc1.Range("D1").Value = 2837.1
$doc1.Range("E1").Value = 10133.9

$doc1.Range("A2").Value = 4781.7
$doc1.Range("B2").Value = 0
$doc1.Range("C2").Value = 2714.6
$doc1.Range("D2").Value = 2861.6
$doc1.Range("E2").Value = 8857.0999999999985

$doc1.Range("A3").Value = 3394.2999999999997
$doc1.Range("B3").Value = 2714.6
$doc1.Range("C3").Value = 0
$doc1.Range("D3").Value = 3045.7
$doc1.Range("E3").Value = 7615.2999999999993

$doc1.Range("A4").Value = 2837.1
$doc1.Range("B4").Value = 2861.6
$doc1.Range("C4").Value = 3045.7
$doc1.Range("D4").Value = 0
$doc1.Range("E4").Value = 10401.299999999999

$doc1.Range("A5").Value = 10133.9
$doc1.Range("B5").Value = 8857.0999999999985
$doc1.Range("C5").Value = 7615.2999999999993
$doc1.Range("D5").Value = 10401.299999999999
$doc1.Range("E5").Value = 0

$doc1.Activate()
$wb.Windows.Item(1).Zoom = 140
$doc1.Range("D10").Select()

# ---------------------------------------------------------------------------
# 4. Add "DOC_2" right after "DOC_1" with its own distance-cost matrix
# ---------------------------------------------------------------------------
$doc1 = $wb.Worksheets.Item("DOC_1")
$doc2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $doc1)
$doc2.Name = "DOC_2"

$doc2.Range("A1").Value = 0
$doc2.Range("B1").Value = 7514.1
$doc2.Range("C1").Value = 5333.9000000000005
$doc2.Range("D1").Value = 4458.3
$doc2.Range("E1").Value = 15924.7

$doc2.Range("A2").Value = 7514.1
$doc2.Range("B2").Value = 0
$doc2.Range("C2").Value = 4265.8
$doc2.Range("D2").Value = 4496.8
$doc2.Range("E2").Value = 13918.300000000001

$doc2.Range("A3").Value = 5333.9000000000005
$doc2.Range("B3").Value = 4265.8
$doc2.Range("C3").Value = 0
$doc2.Range("D3").Value = 4786.1000000000004
$doc2.Range("E3").Value = 11966.900000000001

$doc2.Range("A4").Value = 4458.3
$doc2.Range("B4").Value = 4496.8
$doc2.Range("C4").Value = 4786.1000000000004
$doc2.Range("D4").Value = 0
$doc2.Range("E4").Value = 16344.900000000001

$doc2.Range("A5").Value = 15924.7
$doc2.Range("B5").Value = 13918.300000000001
$doc2.Range("C5").Value = 11966.900000000001
$doc2.Range("D5").Value = 16344.900000000001
$doc2.Range("E5").Value = 0

$doc2.PageSetup.LeftMargin = 0.511811024 * 72
$doc2.PageSetup.RightMargin = 0.511811024 * 72
$doc2.PageSetup.TopMargin = 0.78740157499999996 * 72
$doc2.PageSetup.BottomMargin = 0.78740157499999996 * 72
$doc2.PageSetup.HeaderMargin = 0.31496062000000002 * 72
$doc2.PageSetup.FooterMargin = 0.31496062000000002 * 72

$doc2.Activate()
$wb.Windows.Item(1).Zoom = 130
$doc2.Range("K7").Select()

# ---------------------------------------------------------------------------
# 5. Add "DOC_3" right after "DOC_2" with its own distance-cost matrix
# ---------------------------------------------------------------------------
$doc2 = $wb.Worksheets.Item("DOC_2")
$doc3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $doc2)
$doc3.Name = "DOC_3"

$doc3.Range("A1").Value = 0
$doc3.Range("B1").Value = 9563
$doc3.Range("C1").Value = 6788
$doc3.Range("D1").Value = 5674
$doc3.Range("E1").Value = 20267

$doc3.Range("A2").Value = 9563
$doc3.Range("B2").Value = 0
$doc3.Range("C2").Value = 5429
$doc3.Range("D2").Value = 5723
$doc3.Range("E2").Value = 17714

$doc3.Range("A3").Value = 6788
$doc3.Range("B3").Value = 5429
$doc3.Range("C3").Value = 0
$doc3.Range("D3").Value = 6091
$doc3.Range("E3").Value = 15230

$doc3.Range("A4").Value = 5674
$doc3.Range("B4").Value = 5723
$doc3.Range("C4").Value = 6091
$doc3.Range("D4").Value = 0
$doc3.Range("E4").Value = 20802

$doc3.Range("A5").Value = 20267
$doc3.Range("B5").Value = 17714
$doc3.Range("C5").Value = 15230
$doc3.Range("D5").Value = 20802
$doc3.Range("E5").Value = 0

$doc3.PageSetup.LeftMargin = 0.511811024 * 72
$doc3.PageSetup.RightMargin = 0.511811024 * 72
$doc3.PageSetup.TopMargin = 0.78740157499999996 * 72
$doc3.PageSetup.BottomMargin = 0.78740157499999996 * 72
$doc3.PageSetup.HeaderMargin = 0.31496062000000002 * 72
$doc3.PageSetup.FooterMargin = 0.31496062000000002 * 72

$doc3.Activate()
$wb.Windows.Item(1).Zoom = 120
$doc3.Range("E11").Select()

# ---------------------------------------------------------------------------
# 6. Make "Demands" the active tab and update its selected cell
# ---------------------------------------------------------------------------
$demands = $wb.Worksheets.Item("Demands")
$demands.Activate()
$demands.Range("F12").Select()

foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
